$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 117, shifting rows 117:128 down to 120:131.
$ws.Rows("117:119").Insert()

# Fill the 3 newly inserted rows with the new "Kakamas" variety entries.
# Columns A,B,C,E,F,G,H,I,J,Q,R are constant for this product block (same as
# the existing rows in this "$/bins (400 kilos)" sub-block, e.g. row 116).
$commonRows = 117,118,119
foreach ($r in $commonRows) {
    $ws.Cells.Item($r, 1).Value = 2
    $ws.Cells.Item($r, 2).Value = "Comercializadora del Agro de Limarí"
    $ws.Cells.Item($r, 3).Value = "Coquimbo"
    $ws.Cells.Item($r, 5).Value = 4
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100103
    $ws.Cells.Item($r, 8).Value = "Frutos de hueso (carozo)"
    $ws.Cells.Item($r, 9).Value = 100103004
    $ws.Cells.Item($r, 10).Value = "Durazno"
    $ws.Cells.Item($r, 17).Value = "`$/bins (400 kilos)"
    $ws.Cells.Item($r, 18).Value = "Región de O'Higgins"
}

# Row 117: Kakamas / Especial
$ws.Cells.Item(117, 4).Value = 44644
$ws.Cells.Item(117, 11).Value = "Kakamas"
$ws.Cells.Item(117, 12).Value = "Especial"
$ws.Cells.Item(117, 13).Value = 16
$ws.Cells.Item(117, 14).Value = 470000
$ws.Cells.Item(117, 15).Value = 480000
$ws.Cells.Item(117, 16).Value = 475000
$ws.Cells.Item(117, 19).Value = 1188
$ws.Cells.Item(117, 20).Value = 400

# Row 118: Kakamas / Primera
$ws.Cells.Item(118, 4).Value = 44644
$ws.Cells.Item(118, 11).Value = "Kakamas"
$ws.Cells.Item(118, 12).Value = "Primera"
$ws.Cells.Item(118, 13).Value = 16
$ws.Cells.Item(118, 14).Value = 420000
$ws.Cells.Item(118, 15).Value = 430000
$ws.Cells.Item(118, 16).Value = 425000
$ws.Cells.Item(118, 19).Value = 1062
$ws.Cells.Item(118, 20).Value = 400

# Row 119: Kakamas / Segunda
$ws.Cells.Item(119, 4).Value = 44644
$ws.Cells.Item(119, 11).Value = "Kakamas"
$ws.Cells.Item(119, 12).Value = "Segunda"
$ws.Cells.Item(119, 13).Value = 16
$ws.Cells.Item(119, 14).Value = 400000
$ws.Cells.Item(119, 15).Value = 410000
$ws.Cells.Item(119, 16).Value = 405000
$ws.Cells.Item(119, 19).Value = 1012
$ws.Cells.Item(119, 20).Value = 400
